$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Tear down the old "screen / hex" sizing table (Table1, A1:H2) along with
#    everything it fed (the W/H scratch formulas in I6:L7, the hex-grid
#    offset formulas in B9:C20, and the old H23:H25 scratch column).
# ---------------------------------------------------------------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Delete()

$ws.Range("D1:H2").ClearContents()
$ws.Range("B5:C7").ClearContents()
$ws.Range("I6:L7").ClearContents()
$ws.Range("B9:C12").ClearContents()
$ws.Range("B14:C14").ClearContents()
$ws.Range("B15").ClearContents()
$ws.Range("B17:B18").ClearContents()
$ws.Range("B20:C20").ClearContents()
$ws.Range("H23:H25").ClearContents()

# ---------------------------------------------------------------------------
# 2. Write the new Godot-style "MoveX / MoveY / OffsetX" table data.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value2 = "MoveX"
$ws.Range("B1").Value2 = "MoveY"
$ws.Range("C1").Value2 = "OffsetX"
$ws.Range("A2").Value2 = 198
$ws.Range("B2").Value2 = 170

$newtbl = $ws.ListObjects.Add(1, $ws.Range("A1:C2"), 0, 1)
$newtbl.Name = "Table1"
$ws.Range("C2").Formula = "=Table1[[#This Row],[MoveX]]/2"

# New scratch formula next to the hex grid rows.
$ws.Range("J11").Formula = "=1920/2"

# New explicit zero offsets for the 2,2 hex row.
$ws.Range("B16").Value2 = 0
$ws.Range("C16").Value2 = 0

# ---------------------------------------------------------------------------
# 3. Restore the current selection to match the saved view.
# ---------------------------------------------------------------------------
$ws.Range("B15").Select()
